# Weekly refresh of fruit/vegetable price data:
# a new daily record is inserted at row 65 (shifting all subsequent
# rows, 65-137, down by one, to 66-138), and the sheet's dimension
# grows from A1:R137 to A1:R138.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65; Excel automatically shifts
# every row at/after 65 down by one (old row 65 -> 66, ..., old row
# 137 -> 138) and extends the sheet dimension accordingly.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new data point.
$ws.Range("A65").Value = 8
$ws.Range("B65").Value = "Terminal La Palmera de La Serena"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 44895
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = 100112052
$ws.Range("G65").Value = "Albahaca"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 600
$ws.Range("K65").Value = 4000
$ws.Range("L65").Value = 4500
$ws.Range("M65").Value = 4250
$ws.Range("N65").Value = "$/paquete"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 4250
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"
